$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '27.139.43'
Set-TextValue $ws 'E2' '  -1.16%  '

# Row 3
Set-TextValue $ws 'D3' '1.783.45'
Set-TextValue $ws 'E3' '  -1.80%  '

# Row 4
Set-TextValue $ws 'E4' '  +0.22%  '

# Row 5
Set-TextValue $ws 'D5' '336.80'

# Row 6
Set-TextValue $ws 'E6' '  +0.16%  '

# Row 7
Set-TextValue $ws 'D7' '0.3829'
Set-TextValue $ws 'E7' '  +0.31%  '

# Row 8
Set-TextValue $ws 'D8' '0.3427'
Set-TextValue $ws 'E8' '  -2.18%  '

# Row 9
Set-TextValue $ws 'D9' '47.92'
Set-TextValue $ws 'E9' '  -2.24%  '

# Row 10
Set-TextValue $ws 'D10' '1.192'
Set-TextValue $ws 'E10' '  -3.38%  '

# Row 11
Set-TextValue $ws 'D11' '0.07463'
Set-TextValue $ws 'E11' '  -3.69%  '

# Row 12
Set-TextValue $ws 'D12' '1.001'
Set-TextValue $ws 'E12' '  +0.08%  '

# Row 13
Set-TextValue $ws 'D13' '21.70'
Set-TextValue $ws 'E13' '  -2.69%  '

# Row 14
Set-TextValue $ws 'D14' '6.433'
Set-TextValue $ws 'E14' '  -2.62%  '

# Row 15
Set-TextValue $ws 'D15' '1.781.20'
Set-TextValue $ws 'E15' '  -1.92%  '

# Row 16
Set-TextValue $ws 'D16' '7.112'
Set-TextValue $ws 'E16' '  -1.65%  '

# Row 17
Set-TextValue $ws 'E17' '  -2.32%  '

# Row 18
Set-TextValue $ws 'D18' '0.06648'
Set-TextValue $ws 'E18' '  -1.06%  '

# Row 19
Set-TextValue $ws 'D19' '83.37'
Set-TextValue $ws 'E19' '  -3.23%  '

# Row 20
Set-TextValue $ws 'E20' '  +0.16%  '

# Row 21
Set-TextValue $ws 'D21' '17.44'
Set-TextValue $ws 'E21' '  -0.96%  '

# Row 22
Set-TextValue $ws 'D22' '6.529'
Set-TextValue $ws 'E22' '  -0.67%  '

# Row 23
Set-TextValue $ws 'D23' '27.148.52'
Set-TextValue $ws 'E23' '  -1.12%  '

# Row 24
Set-TextValue $ws 'D24' '12.30'
Set-TextValue $ws 'E24' '  -6.96%  '

# Row 25
Set-TextValue $ws 'D25' '2.366'
Set-TextValue $ws 'E25' '  -4.15%  '

# Row 26
Set-TextValue $ws 'D26' '2.508'
Set-TextValue $ws 'E26' '  -6.06%  '

# Row 27
Set-TextValue $ws 'D27' '21.18'
Set-TextValue $ws 'E27' '  -3.94%  '

# Row 28
Set-TextValue $ws 'E28' '  -1.97%  '

# Row 29
Set-TextValue $ws 'D29' '155.38'
Set-TextValue $ws 'E29' '  +0.84%  '

# Row 30
Set-TextValue $ws 'D30' '1.984.19'
Set-TextValue $ws 'E30' '  -1.77%  '

# Row 31
Set-TextValue $ws 'D31' '134.33'
Set-TextValue $ws 'E31' '  -1.22%  '

# Row 32
Set-TextValue $ws 'D32' '3.976'
Set-TextValue $ws 'E32' '  -1.92%  '

# Row 33
Set-TextValue $ws 'D33' '6.021'
Set-TextValue $ws 'E33' '  -5.10%  '

# Row 34
Set-TextValue $ws 'E34' '  -1.41%  '

# Row 35
Set-TextValue $ws 'E35' '  -6.25%  '

# Row 36
Set-TextValue $ws 'E36' '  -4.11%  '

# Row 37
Set-TextValue $ws 'D37' '5.397'
Set-TextValue $ws 'E37' '  -3.86%  '

# Row 38
Set-TextValue $ws 'D38' '0.6843'
Set-TextValue $ws 'E38' '  -1.88%  '

# Row 39
Set-TextValue $ws 'D39' '0.06330'
Set-TextValue $ws 'E39' '  -2.29%  '

# Row 40
Set-TextValue $ws 'D40' '0.02340'
Set-TextValue $ws 'E40' '  -2.56%  '

# Row 41
Set-TextValue $ws 'D41' '0.2185'
Set-TextValue $ws 'E41' '  -3.53%  '

# Row 42
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D42' '8.430'
Set-TextValue $ws 'E42' '  -5.92%  '

# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D43' '1.237'
Set-TextValue $ws 'E43' '  -4.56%  '

# Row 44
Set-TextValue $ws 'D44' '14.26'
Set-TextValue $ws 'E44' '  -3.88%  '

# Row 45
Set-TextValue $ws 'E45' '  +0.18%  '

# Row 46
Set-TextValue $ws 'D46' '0.6419'
Set-TextValue $ws 'E46' '  -1.84%  '

# Row 47
Set-TextValue $ws 'E47' '  -3.88%  '

# Row 48
Set-TextValue $ws 'D48' '2.159'
Set-TextValue $ws 'E48' '  -0.79%  '

# Row 49
Set-TextValue $ws 'D49' '131.12'
Set-TextValue $ws 'E49' '  -1.32%  '

# Row 50
Set-TextValue $ws 'D50' '0.07106'

# Row 51
Set-TextValue $ws 'D51' '79.00'
Set-TextValue $ws 'E51' '  -1.85%  '
